$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "ScreenRecStarted" category to "0_unstated" wherever it appears
# as (part of) a cell's text value.
$ws.Range("G1").Value = "0_unstated"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Update the active cell selection to G1
$ws.Range("G1").Select()
